$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update author reference for Nascimento from "Nascimento e Teles ()" to "Nascimento et al. (2020)"
$ws.Range("B5").Value = "Nascimento et al. (2020)"

# Widen column B to fit the new, longer text
$ws.Columns.Item(2).ColumnWidth = 26.6

# Update the active cell selection shown in the sheet view
$ws.Range("G10").Select()
